# Rename the "Total" treatment label (column D, rows 20-25) to "Total DST".
# The shared-string table will drop the now-unused "Total" entry and append
# a new "Total DST" entry automatically on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D20").Value = "Total DST"
$ws.Range("D21").Value = "Total DST"
$ws.Range("D22").Value = "Total DST"
$ws.Range("D23").Value = "Total DST"
$ws.Range("D24").Value = "Total DST"
$ws.Range("D25").Value = "Total DST"

# Update the active selection to match the edited workbook (D18).
$ws.Range("D18").Select()
